$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# Inserts a new paragraph directly after the paragraph that contains $anchorText,
# setting $newText as its content in italics (run-level only, not the paragraph mark).
function Add-ItalicParagraphAfterText($anchorText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find anchor text: $anchorText"
    }
    $endPos = $rng.End
    $insertPoint = $d.Range($endPos, $endPos)
    $insertPoint.InsertParagraphAfter() | Out-Null
    $newParaStart = $endPos + 1
    $afterRng = $d.Range($newParaStart, $newParaStart)
    $afterRng.InsertAfter($newText)
    $textEnd = $newParaStart + $newText.Length
    $narrow = $d.Range($newParaStart, $textEnd)
    $narrow.Font.Italic = $true
}

# 1. Title (Heading3) translation swap
Replace-Text "Technology of Biopolymers" "Biopolymers Technology"

# 2. Activation date
Replace-Text "Ativação: 01/01/2017" "Ativação: 01/01/2025"

# 3. Objetivos paragraph
$obj_old = "Apresentar os conceitos básicos da ciência dos polímeros, incluindo as reações químicas e os principais métodos de caracterização envolvidos na preparação destes materiais. Estes fundamentos serão usados para introduzir os alunos aos polímeros de fontes renováveis, produzidos a partir de unidades monoméricas extraídas da biomassa."
$obj_new = "Apresentar os conceitos básicos da ciência dos polímeros e os principais problemas diretamente relacionados ao seu uso e descarte indiscriminados. Desenvolver o pensamento crítico e apresentar ferramentas alternativas para a produção dos polímeros, bem como para minimizar seus impactos ambientais."
Replace-Text $obj_old $obj_new
Add-ItalicParagraphAfterText $obj_new "To present the basic concepts of polymer science and the main problems directly related to their indiscriminate use and disposal. To develop critical thinking and present alternative tools for polymer production, as well as to minimize their environmental impacts."

# 4. Programa resumido paragraph
$res_old = "Fundamentos sobre a química dos polímeros; Mecanismos de polimerização; Caracterização e propriedades gerais dos polímeros; Monômeros derivados da biomassa e principais polímeros obtidos a partir deles."
$res_new = "Reações de polimerização; Propriedades gerais dos polímeros; Monômeros e polímeros derivados de fontes renováveis. Conceitos básicos de circularidade dos materiais poliméricos."
Replace-Text $res_old $res_new
Add-ItalicParagraphAfterText $res_new "Polymerization reactions; General properties of polymers; Monomers and polymers derived from renewable sources. Basic concepts of circularity of polymeric materials."

# 5. Programa paragraph
$prog_old = "Fundamentos sobre a química dos polímeros: composição e estrutura, nomenclatura, polímeros lineares, ramificados, e entrecruzados, massa molar média, propriedades físicas (comportamentos cristalino e amorfo); Mecanismos de polimerização: poliadição e policondensação; Caracterização e propriedades gerais dos polímeros: espectroscopias de infravermelho, FTIR, e ressonância magnética nuclear, RMN, propriedades mecânicas e térmicas; Introdução aos materiais derivados de fontes renováveis; Rotas não-fósseis para a obtenção de monômeros usuais (etileno, propileno, glicerol e derivados); Monômeros exclusivamente obtidos de fontes renováveis (terpenos e terpenóides, monômeros derivados do breu, monômeros derivados de açúcares, ácidos carboxílicos e aminoácidos, furanos, óleos vegetais e derivados); Estado da arte e projeções futuras para os polímeros derivados de fontes renováveis."
$prog_new = "Fundamentos sobre a química dos polímeros: composição e estrutura, massa molecular média, propriedades físicas (comportamentos cristalino e amorfo, propriedades mecânicas e térmicas). Reações de poliadição e policondensação. Introdução aos materiais derivados de fontes renováveis. Monômeros de fonte renovável (etileno, ácidos carboxílicos, aminas, álcoois, óleos vegetais, CO2, entre outros). Polímeros de fonte renovável (celulose, amido, quitina e quitosana, exopolissacarídeos, polihidroxialcanoatos). Estratégias para fim de vida: conceitos básicos de biodegradação e economia circular."
Replace-Text $prog_old $prog_new
Add-ItalicParagraphAfterText $prog_new "Fundamentals of polymer chemistry: composition and structure, average molecular weight, physical properties (crystalline and amorphous behaviours, mechanical and thermal properties). Polyaddition and polycondensation reactions. Introduction to materials derived from renewable sources. Renewable source monomers (ethylene, carboxylic acids, amines, alcohols, vegetable oils, CO2, among others). Polymers directly extracted from renewable sources (cellulose, starch, chitin and chitosan, exopolysaccharides, polyhydroxyalkanoates). End-of-life strategies: basic concepts of biodegradation and circular economy."

# 6. Avaliação - Método
Replace-Text "Duas provas escritas envolvendo o conteúdo teórico ministrado em sala de aula." "Uma avaliação escrita e um estudo de caso."

# 7. Avaliação - Critério
Replace-Text "A nota final corresponderá à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados." "A nota final corresponderá à média aritmética da nota da prova escrita e da nota do estudo de caso. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto aqueles que tiverem média inferior a 3 estarão reprovados."

# 8. Bibliografia
$bib_old = "Eloisa B. Mano; Introdução a Polímeros, Editora Edgard BlücherLtda, São Paulo, 1999.Elizabete F. Lucas, Bluma G. Soares, Elisabeth E. C. Monteiro; Caracterização de polímeros: determinaçãoo de peso molecular e análise térmica. E-papers Serviços Editoriais Ltda, Rio de Janeiro, 2001.Fred J. Davis; PolymerChemistry: a practical approach. Oxford University Press Inc., New York, 2004.George Odian; Principles of Polymerization.John Wiley and Sons, New Jersey, 2004.Mohamed N. Belgacem, Alessandro Gandini; Monomers, polymers and composites from renewable resources.ElsevierLtda, Amsterdam, 2008."
$bib_new = "1-Eloisa B. Mano; Introdução a Polímeros, Editora Edgard BlücherLtda, São Paulo, 1999. 2-2- Sebastião V. Canevarol; Ciência dos Polímeros. Um Texto Básico Para Tecnólogos e Engenheiros. Artliber; 3ª edição. 3-3- J. P. Greene; Sustainable plastics: environmental assessments of biobased, biodegradable, and reclycled plastics. John Wiley & Sons, New Jersey, United States, 2014."
Replace-Text $bib_old $bib_new
